# PP, photos now working for selected
#
# Adds a second example product-code row (A7) styled like the big bold
# title in A3, and adds a new "photo" column header in M1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: big bold title style (same as A3) with "BSET1000-KALORIK" ---
# Copy A3's formatting (bold 24pt font, vertical-centered) onto A7, then set
# its text and restore the big row height used by row 3.
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value2 = "BSET1000-KALORIK"
$ws.Rows(7).RowHeight = 31.5

# --- New column M: "photo" header ---
$ws.Range("M1").Value2 = "photo"

# --- Update the active selection to L3 ---
$ws.Range("L3").Select() | Out-Null
